$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 209.54167
$ws.Range("I33").Value = 158.85
$ws.Range("K33").Value = 158.85
$ws.Range("M33").Value = 70.15

$ws.Range("H107").Value = 1562.2
$ws.Range("I107").Value = 2005
$ws.Range("J107").Value = 1451.5
$ws.Range("K107").Value = 2005
$ws.Range("L107").Value = 1451.5
$ws.Range("M107").Value = -85
$ws.Range("N107").Value = -5291.5

$ws.Range("H137").Value = 2018.3062
$ws.Range("I137").Value = 1576.2354
$ws.Range("J137").Value = 3020.3333
$ws.Range("K137").Value = 4728.706200000001
$ws.Range("L137").Value = 9060.999899999999
$ws.Range("M137").Value = -2178.706200000001
$ws.Range("N137").Value = -14160.9999

$ws.Range("H138").Value = 3746
$ws.Range("I138").Value = 1724.72
$ws.Range("J138").Value = 5430.4
$ws.Range("K138").Value = 5174.16
$ws.Range("L138").Value = 16291.2
$ws.Range("M138").Value = -34.15999999999985
$ws.Range("N138").Value = -26571.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1218.9231
$ws.Range("I2").Value = 1375.9333
$ws.Range("K2").Value = 1375.9333
$ws.Range("M2").Value = -1262.9333

$ws.Range("H116").Value = 1218.9231
$ws.Range("I116").Value = 1375.9333
$ws.Range("K116").Value = 1375.9333
$ws.Range("M116").Value = 918.0667000000001

$ws.Range("H122").Value = 31252252
$ws.Range("I122").Value = 3006
$ws.Range("J122").Value = 62501500
$ws.Range("K122").Value = 9018
$ws.Range("L122").Value = 187504500
$ws.Range("M122").Value = -6568
$ws.Range("N122").Value = -187509400

$ws.Range("H132").Value = 4514.75
$ws.Range("I132").Value = 1648.8334
$ws.Range("J132").Value = 7953.85
$ws.Range("K132").Value = 4946.5002
$ws.Range("L132").Value = 23861.55
$ws.Range("M132").Value = -2416.5002
$ws.Range("N132").Value = -28921.55

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1218.9231
$ws.Range("I3").Value = 1375.9333
$ws.Range("K3").Value = 1375.9333
$ws.Range("M3").Value = -1261.9333

$ws.Range("H105").Value = 1119166.5
$ws.Range("I105").Value = 1490678.6
$ws.Range("K105").Value = 1490678.6
$ws.Range("M105").Value = -1488931.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2173.6182
$ws.Range("I31").Value = 1637.0256
$ws.Range("J31").Value = 3481.5625
$ws.Range("K31").Value = 1637.0256
$ws.Range("L31").Value = 3481.5625
$ws.Range("M31").Value = -1342.0256
$ws.Range("N31").Value = -4071.5625

$ws.Range("H34").Value = 2173.6182
$ws.Range("I34").Value = 1637.0256
$ws.Range("J34").Value = 3481.5625
$ws.Range("K34").Value = 1637.0256
$ws.Range("L34").Value = 3481.5625
$ws.Range("M34").Value = -1435.0256
$ws.Range("N34").Value = -3885.5625

$ws.Range("H109").Value = 49000
$ws.Range("J109").Value = 49000
$ws.Range("L109").Value = 49000
$ws.Range("N109").Value = -51080

$ws.Range("H122").Value = 24319
$ws.Range("I122").Value = 12000
$ws.Range("J122").Value = 36638
$ws.Range("K122").Value = 36000
$ws.Range("L122").Value = 109914
$ws.Range("M122").Value = -33550
$ws.Range("N122").Value = -114814

$ws.Range("H134").Value = 2554.41
$ws.Range("I134").Value = 1584.7354
$ws.Range("J134").Value = 3775.4814
$ws.Range("K134").Value = 4754.206200000001
$ws.Range("L134").Value = 11326.4442
$ws.Range("M134").Value = -2219.206200000001
$ws.Range("N134").Value = -16396.4442

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1055.2
$ws.Range("J22").Value = 1181.5
$ws.Range("L22").Value = 3544.5
$ws.Range("N22").Value = -3882.5

$ws.Range("H27").Value = 1055.2
$ws.Range("J27").Value = 1181.5
$ws.Range("L27").Value = 3544.5
$ws.Range("N27").Value = -3748.5

$ws.Range("H58").Value = 3117.8572
$ws.Range("J58").Value = 3117.8572
$ws.Range("L58").Value = 9353.5716
$ws.Range("N58").Value = -9609.5716

$ws.Range("H112").Value = 2292.5
$ws.Range("I112").Value = 1250
$ws.Range("J112").Value = 4030
$ws.Range("K112").Value = 3750
$ws.Range("L112").Value = 12090
$ws.Range("M112").Value = -2642
$ws.Range("N112").Value = -14306

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 44400.832
$ws.Range("J141").Value = 44400.832
$ws.Range("L141").Value = 44400.832
$ws.Range("N141").Value = -54760.832

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1120
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 1150
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 1150
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -1526

$ws.Range("H100").Value = 4057.2632
$ws.Range("I100").Value = 2358.5
$ws.Range("J100").Value = 5944.778
$ws.Range("K100").Value = 2358.5
$ws.Range("L100").Value = 5944.778
$ws.Range("M100").Value = -1817.5
$ws.Range("N100").Value = -7026.778

$ws.Range("H132").Value = 2874.9268
$ws.Range("I132").Value = 2959.0483
$ws.Range("J132").Value = 2614.15
$ws.Range("K132").Value = 8877.1449
$ws.Range("L132").Value = 7842.450000000001
$ws.Range("M132").Value = -6347.144899999999
$ws.Range("N132").Value = -12902.45

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 44626
$ws.Range("J27").Value = 44626
$ws.Range("L27").Value = 44626
$ws.Range("N27").Value = -44764

$ws.Range("H115").Value = 50377
$ws.Range("J115").Value = 50377
$ws.Range("L115").Value = 50377
$ws.Range("N115").Value = -53511

$ws.Range("H122").Value = 5050.5
$ws.Range("I122").Value = 3399.8333
$ws.Range("K122").Value = 10199.4999
$ws.Range("M122").Value = -7749.499899999999

$ws.Range("H132").Value = 1259.2
$ws.Range("I132").Value = 604.3333
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 1812.9999
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = 717.0001
$ws.Range("N132").Value = -12560

$ws.Range("H136").Value = 2914.578
$ws.Range("I136").Value = 1240.7377
$ws.Range("J136").Value = 6435.4136
$ws.Range("K136").Value = 3722.2131
$ws.Range("L136").Value = 19306.2408
$ws.Range("M136").Value = -1172.2131
$ws.Range("N136").Value = -24406.2408

$ws.Range("H140").Value = 52137.3
$ws.Range("J140").Value = 52137.3
$ws.Range("L140").Value = 52137.3
$ws.Range("N140").Value = -62497.3

$ws.Range("H141").Value = 45511.445
$ws.Range("J141").Value = 45511.445
$ws.Range("L141").Value = 45511.445
$ws.Range("N141").Value = -55871.445

